$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.507.97'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '1.748.37'
$ws.Range("E3").Value = '  -3.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4436'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3609'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07465'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.097'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.14%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("E13").Value = '  -5.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.028'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.140'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.71%  '
$ws.Range("D16").Value = '1.757.82'
$ws.Range("E16").Value = '  -3.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001059'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06404'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.851'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.36%  '
$ws.Range("D23").Value = '27.563.05'
$ws.Range("E23").Value = '  -2.05%  '
$ws.Range("E24").Value = '  -2.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.090'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").Value = '1.954.02'
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.101'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.076'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.648'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08982'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.522'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02318'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2085'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6333'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05973'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.942'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.206'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.760'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.713'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5871'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06860'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.79%  '
